$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 4
$ws.Range("A4").Value = 2
$ws.Range("B5").Value = 2
